$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.089.57"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.779.79"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.545"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.76"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.69%  "
$ws.Range("E9").Value = "  -1.33%  "
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0946"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").Value = "2.035.60"
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").Value = "1.784.00"
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.83%  "
$ws.Range("D15").Value = "34.075.04"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.620"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("D20").Value = "0.0₃0786"
$ws.Range("E20").Value = "  +1.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.06%  "
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("E24").Value = "  -0.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.114"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.73%  "
$ws.Range("E33").Value = "  +3.44%  "
$ws.Range("E34").Value = "  -3.01%  "
$ws.Range("D35").Value = "1.445.06"
$ws.Range("E35").Value = "  +3.41%  "
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("E37").Value = "  +3.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0190"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.93%  "
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("E40").Value = "  +1.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "80.18"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.72"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.915"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.05%  "
$ws.Range("E45").Value = "  +1.36%  "
$ws.Range("E46").Value = "  +1.71%  "
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").Value = "1.937.43"
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("E49").Value = "  -5.72%  "
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "103.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.77%  "
